$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new observation for 2026/01/17 (Sat) was inserted into the daily log,
# pushing the 2026/12/29 .. 2027/01/05 block down by one row (667 -> 668 .. 709 -> 710).
$ws.Rows.Item(668).Insert()

$ws.Range("A668").Value = "'2026/01/17"
$ws.Range("A668").Style = "Normal"
$ws.Range("B668").Value = "土"
$ws.Range("C668").Value = 17
$ws.Range("D668").Value = 179
